# Edit script: adds "PO Forecast" sheet, renames header columns on the
# existing "Weekly Quantity" and "Monthly Trend" sheets.

$wb = $excel.ActiveWorkbook

$wsWeekly = $wb.Worksheets.Item("Weekly Quantity")
$wsMonthly = $wb.Worksheets.Item("Monthly Trend")

# 1) Rename the "Requested quantity" headers.
$wsWeekly.Range("B1").Value = "Weekly_PO_Qty"
$wsMonthly.Range("B1").Value = "Monthly_PO_Qty"

# 2) Add the new "PO Forecast" worksheet right after "Monthly Trend"
#    (i.e. at the end of the workbook).
$wsForecast = $wb.Worksheets.Add([Type]::Missing, $wsMonthly)
$wsForecast.Name = "PO Forecast"

# Header row.
$wsForecast.Range("A1").Value = "ds"
$wsForecast.Range("B1").Value = "PO_Forecast"
$wsForecast.Range("C1").Value = "yhat_lower"
$wsForecast.Range("D1").Value = "yhat_upper"

# Copy the header style (bold, centered, bordered) from the "Weekly Quantity"
# sheet's own header row so the new sheet reuses the same cell style.
$wsWeekly.Range("A1").Copy()
$wsForecast.Range("A1:D1").PasteSpecial(-4122)

# Forecast data rows (ds, PO_Forecast, yhat_lower, yhat_upper).
$forecastData = @(
    @(44934.99999999999,82,-18.62187837640266,196.5644599439464),
    @(44941.99999999999,81,-34.12679005131184,180.9425529133858),
    @(44948.99999999999,81,-29.76540232417062,190.5121220695041),
    @(44962.99999999999,80,-24.00312915421683,187.0677939250865),
    @(44969.99999999999,80,-19.56845206717437,179.063100897504),
    @(44976.99999999999,79,-32.42223164027428,183.2511884178383),
    @(44983.99999999999,79,-19.33798097176469,187.6059107861934),
    @(44990.99999999999,79,-30.36355785060247,183.6923569671785),
    @(44997.99999999999,78,-24.98408688480644,184.4696133209812),
    @(45004.99999999999,78,-27.29602755694077,185.5549128873014),
    @(45018.99999999999,77,-22.1746084461758,184.3401072374897),
    @(45025.99999999999,77,-22.32216395182304,185.5844744313587),
    @(45081.99999999999,74,-33.99051878961365,177.8241419333916),
    @(45088.99999999999,73,-31.54683764198366,182.6734343942362),
    @(45095.99999999999,73,-25.93219471092218,188.1865933165461),
    @(45102.99999999999,72,-39.20148656902175,176.0498318339527),
    @(45109.99999999999,72,-32.13563150508036,176.4403421473006),
    @(45123.99999999999,71,-42.42424426449844,179.9778386974562),
    @(45130.99999999999,71,-39.29975688572635,166.9896313819132),
    @(45137.99999999999,70,-39.78734326564108,172.8226928544079),
    @(45144.99999999999,70,-29.57810152351693,179.7516859755603),
    @(45151.99999999999,70,-37.69918256931583,172.3442933118599),
    @(45165.99999999999,69,-42.96160629006438,176.2639623929389),
    @(45179.99999999999,68,-41.63566339850984,171.6314236777754),
    @(45193.99999999999,67,-42.68392482843421,167.3701082471693),
    @(45214.99999999999,66,-37.12989353516458,172.2488268890319),
    @(45221.99999999999,66,-35.57932172169648,165.4864553720944),
    @(45242.99999999999,65,-45.84511269411114,173.7017668804176),
    @(45249.99999999999,64,-29.59696024943504,172.2974189787947),
    @(45256.99999999999,64,-43.02420221508138,172.7279944943043),
    @(45270.99999999999,63,-37.04293484879292,170.0404039900576),
    @(45298.99999999999,61,-48.47637019210856,167.1096178904),
    @(45305.99999999999,61,-39.91842613917206,173.6461875211845),
    @(45312.99999999999,61,-45.7428509699856,161.0869186329613),
    @(45326.99999999999,60,-42.97372814815434,157.8348165390764),
    @(45333.99999999999,60,-35.63759316067625,164.592878315948),
    @(45347.99999999999,59,-50.2647864041629,160.5687775395335),
    @(45361.99999999999,58,-47.78975293049954,153.4732617433035),
    @(45375.99999999999,57,-50.5249826538404,160.1912176024158),
    @(45382.99999999999,57,-53.78664832412122,162.5204640927272),
    @(45389.99999999999,56,-44.88581174658804,166.426564653765),
    @(45396.99999999999,56,-49.1105391392249,166.6389314169654),
    @(45466.99999999999,52,-46.18985611558607,151.5683435380245),
    @(45480.99999999999,51,-48.35836201927181,166.0010462658734),
    @(45487.99999999999,51,-50.21571650658974,162.5676331066485),
    @(45494.99999999999,51,-55.82099687332357,156.0334522430584),
    @(45501.99999999999,50,-54.29284138511739,142.7582328903843),
    @(45508.99999999999,50,-59.29412143449603,155.6282524809494),
    @(45515.99999999999,49,-54.31206353422103,151.224314059557),
    @(45522.99999999999,49,-55.17058507405499,155.5898480649929),
    @(45529.99999999999,49,-58.04210274109944,154.8481127215407),
    @(45536.99999999999,48,-55.9710627768272,160.8703633354881),
    @(45543.99999999999,48,-62.38169756459065,150.3669864752525)
)

$row = 2
foreach ($rec in $forecastData) {
    $wsForecast.Cells.Item($row, 1).Value = $rec[0]
    $wsForecast.Cells.Item($row, 2).Value = $rec[1]
    $wsForecast.Cells.Item($row, 3).Value = $rec[2]
    $wsForecast.Cells.Item($row, 4).Value = $rec[3]
    $row = $row + 1
}

# Copy the date number-format from the "Weekly Quantity" sheet's date column
# onto the new sheet's "ds" column so it matches the same date style.
$wsWeekly.Range("A2").Copy()
$wsForecast.Range("A2:A54").PasteSpecial(-4122)

Write-Host "PO Forecast sheet added with $($row - 2) data rows."
